# Insert a new price-record row at row 100 (shifting existing rows 100-192
# down to 101-193) and populate it with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 100; Excel shifts rows 100-192 down to
# 101-193 and copies formatting (including the date number format on column D)
# from the row above.
$ws.Rows.Item(100).Insert()

$ws.Cells.Item(100, 1).Value2  = 4
$ws.Cells.Item(100, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(100, 3).Value2  = "Los Lagos"
$ws.Cells.Item(100, 4).Value2  = 45072
$ws.Cells.Item(100, 5).Value2  = 10
$ws.Cells.Item(100, 6).Value2  = 100112052
$ws.Cells.Item(100, 7).Value2  = "Albahaca"
$ws.Cells.Item(100, 8).Value2  = "Sin especificar"
$ws.Cells.Item(100, 9).Value2  = "Primera"
$ws.Cells.Item(100, 10).Value2 = 90
$ws.Cells.Item(100, 11).Value2 = 5000
$ws.Cells.Item(100, 12).Value2 = 5000
$ws.Cells.Item(100, 13).Value2 = 5000
$ws.Cells.Item(100, 14).Value2 = "`$/paquete"
$ws.Cells.Item(100, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(100, 16).Value2 = 5000
$ws.Cells.Item(100, 17).Value2 = 1
$ws.Cells.Item(100, 18).Value2 = "Hortaliza"
